# Refresh the "cryptos" price/volume table with the latest scrape.
# Price (col D) and Volume(1h) (col E) are stored as plain TEXT in the
# source data (not numbers), so numeric-looking Price values are written
# with a leading apostrophe to force Excel to keep them as text instead
# of auto-converting them to the Number type (which would also lose
# precision, e.g. 0.7128 -> 0.71279999999999999).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.305.92"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "1.875.64"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "`'0.7128"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").Value = "`'242.18"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "`'0.3106"
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("D9").Value = "`'0.07744"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "`'24.90"
$ws.Range("E10").Value = "  -0.80%  "
$ws.Range("D11").Value = "`'0.08529"
$ws.Range("E11").Value = "  +3.38%  "
$ws.Range("D12").Value = "1.882.45"
$ws.Range("E12").Value = "  +1.70%  "
$ws.Range("D13").Value = "`'5.215"
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("D14").Value = "`'0.7099"
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("D15").Value = "`'91.41"
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").Value = "29.305.64"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "`'0.000008189"
$ws.Range("E17").Value = "  +5.13%  "
$ws.Range("D18").Value = "`'6.008"
$ws.Range("E18").Value = "  +2.53%  "
$ws.Range("D19").Value = "`'242.09"
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("D20").Value = "2.135.14"
$ws.Range("E20").Value = "  +1.43%  "
$ws.Range("D21").Value = "`'13.24"
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "`'7.808"
$ws.Range("E23").Value = "  -2.24%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "`'0.1606"
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").Value = "`'162.94"
$ws.Range("E26").Value = "  +0.34%  "
$ws.Range("D27").Value = "`'9.049"
$ws.Range("E27").Value = "  +1.55%  "
$ws.Range("D28").Value = "`'18.51"
$ws.Range("E28").Value = "  +1.23%  "
$ws.Range("D29").Value = "`'1.514"
$ws.Range("E29").Value = "  +1.26%  "
$ws.Range("D30").Value = "`'4.398"
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").Value = "`'4.313"
$ws.Range("E31").Value = "  +2.70%  "
$ws.Range("D32").Value = "`'1.281"
$ws.Range("E32").Value = "  -2.54%  "
$ws.Range("D33").Value = "`'0.05251"
$ws.Range("E33").Value = "  +1.31%  "
$ws.Range("D34").Value = "`'1.931"
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("D35").Value = "`'1.175"
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("D36").Value = "`'0.7457"
$ws.Range("E36").Value = "  +2.71%  "
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("D38").Value = "`'0.01866"
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("D39").Value = "`'2.716"
$ws.Range("E39").Value = "  +1.16%  "
$ws.Range("D40").Value = "1.179.54"
$ws.Range("E40").Value = "  +1.75%  "
$ws.Range("D41").Value = "`'6.388"
$ws.Range("E41").Value = "  +3.94%  "
$ws.Range("D42").Value = "`'72.89"
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("D43").Value = "`'0.8871"
$ws.Range("E43").Value = "  -1.83%  "
$ws.Range("D44").Value = "`'106.37"
$ws.Range("E44").Value = "  +4.91%  "
$ws.Range("D45").Value = "`'0.9999"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "2.031.80"
$ws.Range("E46").Value = "  +1.75%  "
$ws.Range("D47").Value = "`'1.809"
$ws.Range("E47").Value = "  +2.47%  "
$ws.Range("D48").Value = "`'0.5208"
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("E49").Value = "  +1.59%  "
$ws.Range("D50").Value = "`'9.389"
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("D51").Value = "`'0.4310"
$ws.Range("E51").Value = "  +1.15%  "
